# Update the EntityId column (A7:A12) to the new small sequential ids.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 5
$ws.Range("A9").Value = 4
$ws.Range("A10").Value = 3
$ws.Range("A11").Value = 2
$ws.Range("A12").Value = 1

# Rows 10-12 previously reflected a "Reset" confirmation state with older
# timestamps; now they share the same "Confirmed" state + new timestamps
# used by rows 7-9 (POST-only endpoints left after the removal).
$ws.Range("E10").Value = "Confirmed"
$ws.Range("E11").Value = "Confirmed"
$ws.Range("E12").Value = "Confirmed"

$ws.Range("F7").Value = "Mon, 13 Sep 2021 05:25:55 GMT"
$ws.Range("F8").Value = "Mon, 13 Sep 2021 05:25:55 GMT"
$ws.Range("F9").Value = "Mon, 13 Sep 2021 05:25:55 GMT"
$ws.Range("F10").Value = "Mon, 13 Sep 2021 05:15:38 GMT"
$ws.Range("F11").Value = "Mon, 13 Sep 2021 05:15:38 GMT"
$ws.Range("F12").Value = "Mon, 13 Sep 2021 05:15:38 GMT"

$ws.Range("G7").Value = "Mon, 13 Sep 2021 05:25:59 GMT"
$ws.Range("G8").Value = "Mon, 13 Sep 2021 05:25:59 GMT"
$ws.Range("G9").Value = "Mon, 13 Sep 2021 05:25:59 GMT"
$ws.Range("G10").Value = "Mon, 13 Sep 2021 05:15:52 GMT"
$ws.Range("G11").Value = "Mon, 13 Sep 2021 05:15:52 GMT"
$ws.Range("G12").Value = "Mon, 13 Sep 2021 05:15:52 GMT"

$ws.Range("H7").Value = "Mon, 13 Sep 2021 05:26:01 GMT"
$ws.Range("H8").Value = "Mon, 13 Sep 2021 05:26:01 GMT"
$ws.Range("H9").Value = "Mon, 13 Sep 2021 05:26:01 GMT"
$ws.Range("H10").Value = "Mon, 13 Sep 2021 05:15:59 GMT"
$ws.Range("H11").Value = "Mon, 13 Sep 2021 05:15:59 GMT"
$ws.Range("H12").Value = "Mon, 13 Sep 2021 05:15:59 GMT"

# Match the author's final selection (cell A7) left in the sheet.
$ws.Range("A7").Select()
